$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting existing rows (and the old row 2,
# 3, 4 data) down by one (mirrors the diff: a new weekly entry is prepended
# ahead of the three existing ones).
$ws.Rows.Item(2).Insert()

# The engine's row Insert() carries over the header row's bold/centered
# formatting onto the new row; strip that back to the plain/default style
# used by the rest of the data rows before re-applying the one style that
# really belongs here (the date format on column D, copied below).
$ws.Range("A2:R2").ClearFormats()

$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = (Get-Date -Year 2021 -Month 10 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 300000000
$ws.Range("G2").Value = "Espárragos"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1773
$ws.Range("N2").Value = "$/paquete"
$ws.Range("O2").Value = "Provincia de Linares"
$ws.Range("P2").Value = 1773
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"

# Copy the date style (custom date format) from the row below onto the new
# row's date cell, matching the original file's per-column styling.
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

# Row 4 (previously "Verde") becomes "Sin especificar" after the shift,
# matching the target data for that week.
$ws.Range("H4").Value = "Sin especificar"
